$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 694
$ws.Range("I18").Value = 694
$ws.Range("K18").Value = 694
$ws.Range("M18").Value = -410
$ws.Range("H40").Value = 3709.6667
$ws.Range("I40").Value = 2321
$ws.Range("K40").Value = 2321
$ws.Range("M40").Value = -2146
$ws.Range("H99").Value = 10233.167
$ws.Range("I99").Value = 279.8
$ws.Range("K99").Value = 839.4000000000001
$ws.Range("M99").Value = 658.5999999999999
$ws.Range("H112").Value = 1829.129
$ws.Range("J112").Value = 1896.963
$ws.Range("L112").Value = 5690.889
$ws.Range("N112").Value = -7906.889
$ws.Range("H125").Value = 3923.6667
$ws.Range("I125").Value = 1160.75
$ws.Range("K125").Value = 10446.75
$ws.Range("M125").Value = -7986.75
$ws.Range("H135").Value = 17859072
$ws.Range("I135").Value = 18520484
$ws.Range("K135").Value = 166684356
$ws.Range("M135").Value = -166681821
$ws.Range("H137").Value = 2421.7073
$ws.Range("I137").Value = 2254.5588
$ws.Range("J137").Value = 3233.5715
$ws.Range("K137").Value = 6763.676399999999
$ws.Range("L137").Value = 9700.7145
$ws.Range("M137").Value = -4213.676399999999
$ws.Range("N137").Value = -14800.7145
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H97").Value = 1205.1389
$ws.Range("I97").Value = 658.05884
$ws.Range("K97").Value = 658.05884
$ws.Range("M97").Value = -162.05884
$ws.Range("H102").Value = 1836
$ws.Range("I102").Value = 1915.5
$ws.Range("K102").Value = 1915.5
$ws.Range("M102").Value = -293.5
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 2076
$ws.Range("I94").Value = 2251
$ws.Range("J94").Value = 413.5
$ws.Range("K94").Value = 2251
$ws.Range("L94").Value = 413.5
$ws.Range("M94").Value = -1800
$ws.Range("N94").Value = -1315.5
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2141.4062
$ws.Range("I31").Value = 2032.2609
$ws.Range("K31").Value = 2032.2609
$ws.Range("M31").Value = -1737.2609
$ws.Range("H34").Value = 2141.4062
$ws.Range("I34").Value = 2032.2609
$ws.Range("K34").Value = 2032.2609
$ws.Range("M34").Value = -1830.2609
$ws.Range("H51").Value = 52000
$ws.Range("I51").Value = 52000
$ws.Range("J51").Value = 0
$ws.Range("K51").Value = 52000
$ws.Range("L51").Value = 0
$ws.Range("M51").Value = -51264
$ws.Range("N51").ClearContents()
$ws.Range("H58").Value = 2808.4285
$ws.Range("I58").Value = 1691.9231
$ws.Range("K58").Value = 1691.9231
$ws.Range("M58").Value = -1488.9231
$ws.Range("H60").Value = 7745
$ws.Range("I60").Value = 7745
$ws.Range("K60").Value = 7745
$ws.Range("M60").Value = -7234
$ws.Range("H61").Value = 52000
$ws.Range("I61").Value = 52000
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 52000
$ws.Range("L61").Value = 0
$ws.Range("M61").Value = -51652
$ws.Range("N61").ClearContents()
$ws.Range("H86").Value = 72492.5
$ws.Range("I86").Value = 105995.664
$ws.Range("J86").Value = 38989.332
$ws.Range("K86").Value = 105995.664
$ws.Range("L86").Value = 38989.332
$ws.Range("M86").Value = -104872.664
$ws.Range("N86").Value = -41235.332
$ws.Range("H89").Value = 72492.5
$ws.Range("I89").Value = 105995.664
$ws.Range("J89").Value = 38989.332
$ws.Range("K89").Value = 529978.3200000001
$ws.Range("L89").Value = 194946.66
$ws.Range("M89").Value = -524362.3200000001
$ws.Range("N89").Value = -206178.66
$ws.Range("H98").Value = 0
$ws.Range("I98").Value = 0
$ws.Range("J98").Value = 0
$ws.Range("K98").Value = 0
$ws.Range("L98").Value = 0
$ws.Range("M98").ClearContents()
$ws.Range("N98").ClearContents()
$ws.Range("H99").Value = 44888052
$ws.Range("J99").Value = 50012504
$ws.Range("L99").Value = 50012504
$ws.Range("N99").Value = -50015500
$ws.Range("H105").Value = 3172.9285
$ws.Range("I105").Value = 2875.0833
$ws.Range("K105").Value = 2875.0833
$ws.Range("M105").Value = -1128.0833
$ws.Range("H109").Value = 42428.332
$ws.Range("J109").Value = 42428.332
$ws.Range("L109").Value = 42428.332
$ws.Range("N109").Value = -44508.332
$ws.Range("H117").Value = 72500
$ws.Range("J117").Value = 72500
$ws.Range("L117").Value = 72500
$ws.Range("N117").Value = -81678
$ws.Range("H125").Value = 125000
$ws.Range("J125").Value = 125000
$ws.Range("L125").Value = 125000
$ws.Range("N125").Value = -129920
$ws.Range("H126").Value = 44888052
$ws.Range("J126").Value = 50012504
$ws.Range("L126").Value = 150037512
$ws.Range("N126").Value = -150042452
$ws.Range("H136").Value = 2808.4285
$ws.Range("I136").Value = 1691.9231
$ws.Range("K136").Value = 5075.7693
$ws.Range("M136").Value = -2525.7693
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H62").Value = 0
$ws.Range("I62").Value = 0
$ws.Range("K62").Value = 0
$ws.Range("M62").ClearContents()
$ws.Range("H63").Value = 5051.4443
$ws.Range("J63").Value = 5428.4287
$ws.Range("L63").Value = 16285.2861
$ws.Range("N63").Value = -17783.2861
$ws.Range("H64").Value = 1339.6
$ws.Range("I64").Value = 1339.6
$ws.Range("K64").Value = 4018.8
$ws.Range("M64").Value = -3748.8
$ws.Range("H65").Value = 0
$ws.Range("I65").Value = 0
$ws.Range("K65").Value = 0
$ws.Range("M65").ClearContents()
$ws.Range("H66").Value = 5051.4443
$ws.Range("J66").Value = 5428.4287
$ws.Range("L66").Value = 48855.85830000001
$ws.Range("N66").Value = -56343.85830000001
$ws.Range("H67").Value = 1339.6
$ws.Range("I67").Value = 1339.6
$ws.Range("K67").Value = 4018.8
$ws.Range("M67").Value = -3082.8
$ws.Range("H75").Value = 732.4
$ws.Range("I75").Value = 157
$ws.Range("K75").Value = 471
$ws.Range("M75").Value = 527
$ws.Range("H78").Value = 732.4
$ws.Range("I78").Value = 157
$ws.Range("K78").Value = 1413
$ws.Range("M78").Value = 3579
$ws.Range("H140").Value = 22733090
$ws.Range("I140").Value = 62503500
$ws.Range("J140").Value = 7142.857
$ws.Range("K140").Value = 187510500
$ws.Range("L140").Value = 21428.571
$ws.Range("M140").Value = -187505320
$ws.Range("N140").Value = -31788.571
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H12").Value = 3934.6667
$ws.Range("I12").Value = 800
$ws.Range("J12").Value = 5502
$ws.Range("K12").Value = 800
$ws.Range("L12").Value = 5502
$ws.Range("N12").Value = -5782
$ws.Range("M12").Value = -660
$ws.Range("H14").Value = 797833.25
$ws.Range("I14").Value = 1250937.5
$ws.Range("J14").Value = 279999.84
$ws.Range("K14").Value = 1250937.5
$ws.Range("L14").Value = 279999.84
$ws.Range("M14").Value = -1250769.5
$ws.Range("N14").Value = -280335.84
$ws.Range("H32").Value = 60290
$ws.Range("J32").Value = 60290
$ws.Range("L32").Value = 60290
$ws.Range("N32").Value = -60882
$ws.Range("H42").Value = 49215
$ws.Range("J42").Value = 49215
$ws.Range("L42").Value = 49215
$ws.Range("N42").Value = -50185
$ws.Range("H115").Value = 49215
$ws.Range("J115").Value = 49215
$ws.Range("L115").Value = 49215
$ws.Range("N115").Value = -51565
$ws.Range("H132").Value = 2418.1365
$ws.Range("I132").Value = 2041.4706
$ws.Range("J132").Value = 3698.8
$ws.Range("K132").Value = 6124.4118
$ws.Range("L132").Value = 11096.4
$ws.Range("M132").Value = -3594.4118
$ws.Range("N132").Value = -16156.4
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H17").Value = 432.5
$ws.Range("I17").Value = 432.5
$ws.Range("K17").Value = 432.5
$ws.Range("M17").Value = -262.5
$ws.Range("H29").Value = 4222
$ws.Range("J29").Value = 4222
$ws.Range("L29").Value = 4222
$ws.Range("N29").Value = -4812
$ws.Range("H35").Value = 976.2
$ws.Range("I35").Value = 293.66666
$ws.Range("J35").Value = 2000
$ws.Range("K35").Value = 293.66666
$ws.Range("L35").Value = 2000
$ws.Range("M35").Value = 42.33334000000002
$ws.Range("N35").Value = -2672
$ws.Range("H46").Value = 1749.25
$ws.Range("I46").Value = 1001
$ws.Range("K46").Value = 1001
$ws.Range("M46").Value = -813
$ws.Range("H132").Value = 3485.8838
$ws.Range("I132").Value = 2133.9119
$ws.Range("J132").Value = 8593.333000000001
$ws.Range("K132").Value = 6401.7357
$ws.Range("L132").Value = 25779.999
$ws.Range("M132").Value = -3871.7357
$ws.Range("N132").Value = -30839.999
$ws.Range("H136").Value = 3240.4666
$ws.Range("I136").Value = 2757.6428
$ws.Range("J136").Value = 10000
$ws.Range("K136").Value = 8272.928400000001
$ws.Range("L136").Value = 30000
$ws.Range("M136").Value = -5722.928400000001
$ws.Range("N136").Value = -35100
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H70").Value = 43333
$ws.Range("J70").Value = 43333
$ws.Range("L70").Value = 43333
$ws.Range("N70").Value = -43963
$ws.Range("H73").Value = 43333
$ws.Range("J73").Value = 43333
$ws.Range("L73").Value = 43333
$ws.Range("N73").Value = -45517
$ws.Range("H122").Value = 3479.889
$ws.Range("I122").Value = 3566.125
$ws.Range("J122").Value = 2790
$ws.Range("K122").Value = 10698.375
$ws.Range("L122").Value = 8370
$ws.Range("M122").Value = -8248.375
$ws.Range("N122").Value = -13270
